$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed odds cells in existing rows 2-5 ---
# Row 2
$ws.Cells.Item(2, 7).Value = 1.48
$ws.Cells.Item(2, 8).Value = 4.5
$ws.Cells.Item(2, 9).Value = 6.5
$ws.Cells.Item(2, 10).Value = 2
$ws.Cells.Item(2, 11).Value = 2.4
$ws.Cells.Item(2, 13).Value = 1.04
$ws.Cells.Item(2, 14).Value = 13
$ws.Cells.Item(2, 15).Value = 1.2
$ws.Cells.Item(2, 16).Value = 4.33
$ws.Cells.Item(2, 17).Value = 1.7
$ws.Cells.Item(2, 18).Value = 2.1
$ws.Cells.Item(2, 19).Value = 1.33
$ws.Cells.Item(2, 20).Value = 3.25
$ws.Cells.Item(2, 21).Value = 1.83
$ws.Cells.Item(2, 22).Value = 1.83
$ws.Cells.Item(2, 23).Value = 7.5
$ws.Cells.Item(2, 24).Value = 7.5
$ws.Cells.Item(2, 26).Value = 10
$ws.Cells.Item(2, 27).Value = 12
$ws.Cells.Item(2, 28).Value = 23
$ws.Cells.Item(2, 29).Value = 13
$ws.Cells.Item(2, 30).Value = 8.5
$ws.Cells.Item(2, 31).Value = 17
$ws.Cells.Item(2, 33).Value = 251
$ws.Cells.Item(2, 34).Value = 19
$ws.Cells.Item(2, 35).Value = 34
$ws.Cells.Item(2, 36).Value = 21
$ws.Cells.Item(2, 39).Value = 41
$ws.Cells.Item(2, 40).Value = 3.5
$ws.Cells.Item(2, 41).Value = 7
$ws.Cells.Item(2, 42).Value = 17
$ws.Cells.Item(2, 43).Value = 21
$ws.Cells.Item(2, 44).Value = 41
$ws.Cells.Item(2, 45).Value = 126
$ws.Cells.Item(2, 46).Value = 3.25
$ws.Cells.Item(2, 47).Value = 8.5
$ws.Cells.Item(2, 49).Value = 8
$ws.Cells.Item(2, 51).Value = 34
$ws.Cells.Item(2, 53).Value = 126
$ws.Cells.Item(2, 54).Value = 251

# Row 3
$ws.Cells.Item(3, 7).Value = 2.05
$ws.Cells.Item(3, 8).Value = 3.1
$ws.Cells.Item(3, 9).Value = 3.7
$ws.Cells.Item(3, 10).Value = 3
$ws.Cells.Item(3, 11).Value = 1.91
$ws.Cells.Item(3, 12).Value = 4.5
$ws.Cells.Item(3, 13).Value = 1.1
$ws.Cells.Item(3, 14).Value = 7
$ws.Cells.Item(3, 15).Value = 1.5
$ws.Cells.Item(3, 16).Value = 2.5
$ws.Cells.Item(3, 17).Value = 2.5
$ws.Cells.Item(3, 18).Value = 1.5
$ws.Cells.Item(3, 19).Value = 1.57
$ws.Cells.Item(3, 20).Value = 2.25
$ws.Cells.Item(3, 21).Value = 2.2
$ws.Cells.Item(3, 22).Value = 1.62
$ws.Cells.Item(3, 24).Value = 8.5
$ws.Cells.Item(3, 25).Value = 10
$ws.Cells.Item(3, 26).Value = 19
$ws.Cells.Item(3, 27).Value = 21
$ws.Cells.Item(3, 29).Value = 6.5
$ws.Cells.Item(3, 34).Value = 8
$ws.Cells.Item(3, 35).Value = 17
$ws.Cells.Item(3, 36).Value = 15
$ws.Cells.Item(3, 37).Value = 41
$ws.Cells.Item(3, 38).Value = 41
$ws.Cells.Item(3, 40).Value = 4
$ws.Cells.Item(3, 41).Value = 13
$ws.Cells.Item(3, 42).Value = 29
$ws.Cells.Item(3, 43).Value = 41
$ws.Cells.Item(3, 44).Value = 81
$ws.Cells.Item(3, 45).Value = 301
$ws.Cells.Item(3, 46).Value = 2.25
$ws.Cells.Item(3, 47).Value = 9.5
$ws.Cells.Item(3, 49).Value = 5.5
$ws.Cells.Item(3, 50).Value = 23
$ws.Cells.Item(3, 52).Value = 81

# Row 4
$ws.Cells.Item(4, 7).Value = 4.5
$ws.Cells.Item(4, 9).Value = 2
$ws.Cells.Item(4, 10).Value = 5
$ws.Cells.Item(4, 11).Value = 1.91
$ws.Cells.Item(4, 12).Value = 2.75
$ws.Cells.Item(4, 13).Value = 1.13
$ws.Cells.Item(4, 14).Value = 6
$ws.Cells.Item(4, 15).Value = 1.5
$ws.Cells.Item(4, 16).Value = 2.5
$ws.Cells.Item(4, 17).Value = 2.6
$ws.Cells.Item(4, 18).Value = 1.48
$ws.Cells.Item(4, 19).Value = 1.57
$ws.Cells.Item(4, 20).Value = 2.25
$ws.Cells.Item(4, 21).Value = 2.2
$ws.Cells.Item(4, 22).Value = 1.62
$ws.Cells.Item(4, 23).Value = 9
$ws.Cells.Item(4, 24).Value = 19
$ws.Cells.Item(4, 25).Value = 15
$ws.Cells.Item(4, 27).Value = 41
$ws.Cells.Item(4, 28).Value = 51
$ws.Cells.Item(4, 29).Value = 6
$ws.Cells.Item(4, 31).Value = 19
$ws.Cells.Item(4, 32).Value = 81
$ws.Cells.Item(4, 34).Value = 5.5
$ws.Cells.Item(4, 35).Value = 8
$ws.Cells.Item(4, 36).Value = 9.5
$ws.Cells.Item(4, 37).Value = 17
$ws.Cells.Item(4, 39).Value = 41
$ws.Cells.Item(4, 40).Value = 6
$ws.Cells.Item(4, 41).Value = 26
$ws.Cells.Item(4, 42).Value = 41
$ws.Cells.Item(4, 43).Value = 101
$ws.Cells.Item(4, 44).Value = 151
$ws.Cells.Item(4, 45).Value = 351
$ws.Cells.Item(4, 46).Value = 2.25
$ws.Cells.Item(4, 47).Value = 9.5
$ws.Cells.Item(4, 48).Value = 81
$ws.Cells.Item(4, 49).Value = 3.75
$ws.Cells.Item(4, 50).Value = 12

# Row 5
$ws.Cells.Item(5, 7).Value = 3.6
$ws.Cells.Item(5, 8).Value = 3.1
$ws.Cells.Item(5, 9).Value = 2.15
$ws.Cells.Item(5, 10).Value = 4.33
$ws.Cells.Item(5, 11).Value = 1.95
$ws.Cells.Item(5, 12).Value = 3
$ws.Cells.Item(5, 13).Value = 1.1
$ws.Cells.Item(5, 14).Value = 7
$ws.Cells.Item(5, 15).Value = 1.44
$ws.Cells.Item(5, 16).Value = 2.63
$ws.Cells.Item(5, 17).Value = 2.4
$ws.Cells.Item(5, 18).Value = 1.53
$ws.Cells.Item(5, 19).Value = 1.57
$ws.Cells.Item(5, 20).Value = 2.25
$ws.Cells.Item(5, 21).Value = 2.1
$ws.Cells.Item(5, 22).Value = 1.67
$ws.Cells.Item(5, 23).Value = 8
$ws.Cells.Item(5, 24).Value = 17
$ws.Cells.Item(5, 25).Value = 13
$ws.Cells.Item(5, 26).Value = 41
$ws.Cells.Item(5, 27).Value = 34
$ws.Cells.Item(5, 29).Value = 6.5
$ws.Cells.Item(5, 31).Value = 19
$ws.Cells.Item(5, 32).Value = 67
$ws.Cells.Item(5, 33).Value = 900
$ws.Cells.Item(5, 34).Value = 6
$ws.Cells.Item(5, 35).Value = 9
$ws.Cells.Item(5, 36).Value = 10
$ws.Cells.Item(5, 37).Value = 19
$ws.Cells.Item(5, 39).Value = 41
$ws.Cells.Item(5, 40).Value = 5.5
$ws.Cells.Item(5, 41).Value = 21
$ws.Cells.Item(5, 42).Value = 34
$ws.Cells.Item(5, 43).Value = 81
$ws.Cells.Item(5, 44).Value = 126
$ws.Cells.Item(5, 45).Value = 351
$ws.Cells.Item(5, 46).Value = 2.25
$ws.Cells.Item(5, 47).Value = 9
$ws.Cells.Item(5, 49).Value = 4
$ws.Cells.Item(5, 51).Value = 29
$ws.Cells.Item(5, 53).Value = 81
$ws.Cells.Item(5, 54).Value = 251

# --- Append new rows 6 and 7 ---
# Row 6
$ws.Cells.Item(6, 2).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "xdKSpQD6"
$ws.Cells.Item(6, 2).Value = "12/11/2024"
$ws.Cells.Item(6, 3).Value = "19:30"
$ws.Cells.Item(6, 4).Value = "VENEZUELA - LIGA FUTVE"
$ws.Cells.Item(6, 5).Value = "Caracas"
$ws.Cells.Item(6, 6).Value = "Dep. Tachira"
$ws.Cells.Item(6, 7).Value = 3
$ws.Cells.Item(6, 8).Value = 2.7
$ws.Cells.Item(6, 9).Value = 2.57
$ws.Cells.Item(6, 10).Value = 3.6
$ws.Cells.Item(6, 11).Value = 1.85
$ws.Cells.Item(6, 12).Value = 3.3
$ws.Cells.Item(6, 13).Value = 1.08
$ws.Cells.Item(6, 14).Value = 5.25
$ws.Cells.Item(6, 15).Value = 1.5
$ws.Cells.Item(6, 16).Value = 2.25
$ws.Cells.Item(6, 17).Value = 2.42
$ws.Cells.Item(6, 18).Value = 1.44
$ws.Cells.Item(6, 19).Value = 1.53
$ws.Cells.Item(6, 20).Value = 2.18
$ws.Cells.Item(6, 21).Value = 1.98
$ws.Cells.Item(6, 22).Value = 1.65
$ws.Cells.Item(6, 23).Value = 7.2
$ws.Cells.Item(6, 24).Value = 14.5
$ws.Cells.Item(6, 25).Value = 11
$ws.Cells.Item(6, 26).Value = 40
$ws.Cells.Item(6, 27).Value = 32
$ws.Cells.Item(6, 28).Value = 45
$ws.Cells.Item(6, 29).Value = 6
$ws.Cells.Item(6, 30).Value = 5.4
$ws.Cells.Item(6, 31).Value = 16.5
$ws.Cells.Item(6, 32).Value = 100
$ws.Cells.Item(6, 33).Value = 201
$ws.Cells.Item(6, 34).Value = 6.1
$ws.Cells.Item(6, 35).Value = 11.25
$ws.Cells.Item(6, 36).Value = 10.25
$ws.Cells.Item(6, 37).Value = 30
$ws.Cells.Item(6, 38).Value = 28
$ws.Cells.Item(6, 39).Value = 45
$ws.Cells.Item(6, 40).Value = 4.65
$ws.Cells.Item(6, 41).Value = 17
$ws.Cells.Item(6, 42).Value = 27
$ws.Cells.Item(6, 43).Value = 90
$ws.Cells.Item(6, 44).Value = 150
$ws.Cells.Item(6, 45).Value = 400
$ws.Cells.Item(6, 46).Value = 2.15
$ws.Cells.Item(6, 47).Value = 7.3
$ws.Cells.Item(6, 48).Value = 80
$ws.Cells.Item(6, 49).Value = 4.3
$ws.Cells.Item(6, 50).Value = 15
$ws.Cells.Item(6, 51).Value = 27
$ws.Cells.Item(6, 52).Value = 75
$ws.Cells.Item(6, 53).Value = 150
$ws.Cells.Item(6, 54).Value = 450
$ws.Cells.Item(6, 55).Value = 51
$ws.Cells.Item(6, 56).Value = 51

# Row 7
$ws.Cells.Item(7, 2).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "bB8yq4rJ"
$ws.Cells.Item(7, 2).Value = "12/11/2024"
$ws.Cells.Item(7, 3).Value = "19:30"
$ws.Cells.Item(7, 4).Value = "VENEZUELA - LIGA FUTVE"
$ws.Cells.Item(7, 5).Value = "Zamora"
$ws.Cells.Item(7, 6).Value = "Rayo Zuliano"
$ws.Cells.Item(7, 7).Value = 2.57
$ws.Cells.Item(7, 8).Value = 3.1
$ws.Cells.Item(7, 9).Value = 2.65
$ws.Cells.Item(7, 10).Value = 3.15
$ws.Cells.Item(7, 11).Value = 2.05
$ws.Cells.Item(7, 12).Value = 3.2
$ws.Cells.Item(7, 13).Value = 1.01
$ws.Cells.Item(7, 14).Value = 7.9
$ws.Cells.Item(7, 15).Value = 1.35
$ws.Cells.Item(7, 16).Value = 2.7
$ws.Cells.Item(7, 17).Value = 2.02
$ws.Cells.Item(7, 18).Value = 1.62
$ws.Cells.Item(7, 19).Value = 1.4
$ws.Cells.Item(7, 20).Value = 2.52
$ws.Cells.Item(7, 21).Value = 1.78
$ws.Cells.Item(7, 22).Value = 1.83
$ws.Cells.Item(7, 23).Value = 7.5
$ws.Cells.Item(7, 24).Value = 12.5
$ws.Cells.Item(7, 25).Value = 9.75
$ws.Cells.Item(7, 26).Value = 29
$ws.Cells.Item(7, 27).Value = 23
$ws.Cells.Item(7, 28).Value = 35
$ws.Cells.Item(7, 29).Value = 8.25
$ws.Cells.Item(7, 30).Value = 6
$ws.Cells.Item(7, 31).Value = 14.5
$ws.Cells.Item(7, 32).Value = 75
$ws.Cells.Item(7, 33).Value = 600
$ws.Cells.Item(7, 34).Value = 7.9
$ws.Cells.Item(7, 35).Value = 13
$ws.Cells.Item(7, 36).Value = 10
$ws.Cells.Item(7, 37).Value = 30
$ws.Cells.Item(7, 38).Value = 23
$ws.Cells.Item(7, 39).Value = 35
$ws.Cells.Item(7, 40).Value = 4.45
$ws.Cells.Item(7, 41).Value = 14
$ws.Cells.Item(7, 42).Value = 22
$ws.Cells.Item(7, 43).Value = 60
$ws.Cells.Item(7, 44).Value = 100
$ws.Cells.Item(7, 45).Value = 300
$ws.Cells.Item(7, 46).Value = 2.47
$ws.Cells.Item(7, 47).Value = 6.9
$ws.Cells.Item(7, 48).Value = 60
$ws.Cells.Item(7, 49).Value = 4.5
$ws.Cells.Item(7, 50).Value = 14
$ws.Cells.Item(7, 51).Value = 21
$ws.Cells.Item(7, 52).Value = 60
$ws.Cells.Item(7, 53).Value = 90
$ws.Cells.Item(7, 54).Value = 250
$ws.Cells.Item(7, 55).Value = 51
$ws.Cells.Item(7, 56).Value = 51

